# Update cryptocurrency price (D) and 1h volume change (E) columns
# with freshly scraped values. Cells are plain text (inlineStr in the
# original workbook), so we force a Text number format before writing
# the value - this stops Excel from re-interpreting strings such as
# "0.671" or "43.711.48" as numbers/dates - and then restore the
# default "Normal" style so no stray cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "43.711.48"
Set-TextValue "E2" "  +0.09%  "

Set-TextValue "D3" "2.348.05"
Set-TextValue "E3" "  +0.26%  "

Set-TextValue "E4" "  -0.05%  "

Set-TextValue "D5" "0.671"
Set-TextValue "E5" "  +3.26%  "

Set-TextValue "D6" "234.73"
Set-TextValue "E6" "  +0.39%  "

Set-TextValue "D7" "73.51"
Set-TextValue "E7" "  +9.19%  "

Set-TextValue "E8" "  -0.05%  "

Set-TextValue "D9" "0.541"
Set-TextValue "E9" "  +17.82%  "

Set-TextValue "D10" "0.0980"
Set-TextValue "E10" "  +0.75%  "

Set-TextValue "D11" "28.26"
Set-TextValue "E11" "  +3.84%  "

Set-TextValue "E12" "  +1.65%  "

Set-TextValue "D13" "2.698.63"
Set-TextValue "E13" "  +0.37%  "

Set-TextValue "D14" "16.63"
Set-TextValue "E14" "  +6.57%  "

Set-TextValue "D15" "6.68"
Set-TextValue "E15" "  +6.77%  "

Set-TextValue "D16" "0.887"
Set-TextValue "E16" "  +3.97%  "

Set-TextValue "D17" "2.360.86"
Set-TextValue "E17" "  +0.67%  "

Set-TextValue "D18" "43.743.08"
Set-TextValue "E18" "  +0.30%  "

Set-TextValue "E19" "  +2.95%  "

Set-TextValue "D20" "76.82"
Set-TextValue "E20" "  +3.31%  "

Set-TextValue "D21" "6.38"
Set-TextValue "E21" "  +1.46%  "

Set-TextValue "D22" "252.55"
Set-TextValue "E22" "  +1.04%  "

Set-TextValue "D24" "3.75"
Set-TextValue "E24" "  -1.53%  "

Set-TextValue "D25" "2.47"
Set-TextValue "E25" "  +1.46%  "

Set-TextValue "D26" "10.54"
Set-TextValue "E26" "  +5.31%  "

Set-TextValue "E27" "  +0.82%  "

Set-TextValue "D28" "22.31"
Set-TextValue "E28" "  +0.27%  "

Set-TextValue "E29" "  +7.72%  "

Set-TextValue "D30" "172.60"
Set-TextValue "E30" "  -1.22%  "

Set-TextValue "D31" "0.130"
Set-TextValue "E31" "  +0.20%  "

Set-TextValue "E32" "  +4.67%  "

Set-TextValue "E33" "  +2.58%  "

Set-TextValue "E34" "  +3.23%  "

Set-TextValue "D35" "5.15"
Set-TextValue "E35" "  +3.16%  "

Set-TextValue "D36" "3.84"
Set-TextValue "E36" "  +7.01%  "

Set-TextValue "E37" "  -5.74%  "

Set-TextValue "D38" "6.37"
Set-TextValue "E38" "  -3.05%  "

Set-TextValue "E39" "  +5.48%  "

Set-TextValue "D40" "19.38"
Set-TextValue "E40" "  +5.47%  "

Set-TextValue "E41" "  -0.06%  "

Set-TextValue "D42" "8.85"
Set-TextValue "E42" "  -2.38%  "

Set-TextValue "D43" "0.0975"
Set-TextValue "E43" "  +2.52%  "

Set-TextValue "E44" "  +1.53%  "

Set-TextValue "E45" "  -2.30%  "

Set-TextValue "E46" "  +11.87%  "

Set-TextValue "E47" "  +1.73%  "

Set-TextValue "D48" "96.93"
Set-TextValue "E48" "  -2.81%  "

Set-TextValue "D49" "1.431.27"
Set-TextValue "E49" "  -1.22%  "

Set-TextValue "D51" "2.569.34"
Set-TextValue "E51" "  +0.30%  "
